# "added quirks and languages choice option"
#
# odmieniecAncestry sheet (sheet1): random-roll attribute values were
# rerolled, the extra "Elifcki" language note next to "Język Mowa:" was
# removed, and "Rozmiar:" now stores a plain numeric 1 instead of the
# text "0.5". The previously-active tab (Pozorne Pochodzenie) is swapped
# back to the main odmieniecAncestry sheet, with a new selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("odmieniecAncestry")

$ws.Range("B2").Value  = 9    # Sila:
$ws.Range("B3").Value  = 10   # Zrecznosc:
$ws.Range("B5").Value  = 10   # Wola:
$ws.Range("B6").Value  = 11   # Percepcja:
$ws.Range("B7").Value  = 9    # Zdrowie:
$ws.Range("B9").Value  = 10   # Obrona:

$ws.Range("C15").ClearContents()   # drop the extra "Elifcki" language note
$ws.Range("B17").Value = 1         # Rozmiar: text "0.5" -> number 1

$ws.Activate()
$ws.Range("F41").Select()
